# Auto-generated edit script applying the diff to before.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "'43.9583"
$ws.Range("C2").Value = "'-80.0473"
$ws.Range("D2").Value = "43.9583, -80.0473"
$ws.Range("F2").Value = "The Proponent must implement vegetation-enhanced stone revetment and catchbasin shields as offsetting measures by September 30, 2021, with monitoring for vegetation survival (80%) and turbidity reduction (50%)."
$ws.Range("G2").Value = "The Proponent must not carry out any activities that adversely impact the offsetting measures and must provide access permissions for DFO to monitor them."
$ws.Range("H2").Value = "The authorization may be revoked or amended if necessary to protect species at risk, and compliance is mandatory to avoid legal consequences."
$ws.Range("I2").Value = "The Proponent must submit annual monitoring reports for two years post-construction and ensure all offset measures meet effectiveness criteria."
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = "Riparian vegetation (e.g. trees and shrubs and grass)"
$ws.Range("N2").Value = "None"
$ws.Range("O2").Value = "None"
$ws.Range("P2").Value = "Riffles (rivers)"
# Row 3
$ws.Range("F3").Value = "The Proponent must implement offsetting measures to address serious harm to fish, including habitat restoration and creation as per the operating plan."
$ws.Range("G3").Value = "Annual monitoring reports required for dam operations, fish species, and contingency measures, with specific submission deadlines and data formats."
$ws.Range("H3").Value = "Prohibition on depositing deleterious substances and compliance with Species at Risk Act requirements."
$ws.Range("I3").Value = "Authorization cannot be transferred without prior notification to DFO, and non-compliance may result in legal action."
$ws.Range("P3").Value = "Riffles (rivers)"
# Row 4
$ws.Range("E4").Value = "MAY 08 2015"
$ws.Range("F4").Value = "The Proponent must complete offsetting measures during construction phase, including habitat enhancements and rock shoal installation."
$ws.Range("G4").Value = "Monitoring and reporting requirements include annual reports by March 1st, detailing fish rescue efforts, habitat effectiveness, and any adaptive measures."
$ws.Range("H4").Value = "No adverse impacts on offsetting measures are allowed, and compliance with other regulations (e.g., SARA) is mandatory."
$ws.Range("I4").Value = "Authorization cannot be transferred without prior notification to DFO, and all work must adhere to design and safety standards."
$ws.Range("J4").Value = "Riverine"
$ws.Range("K4").Value = "Walleye"
$ws.Range("L4").Value = 2095
$ws.Range("N4").Value = "None"
$ws.Range("P4").Value = "Rock shoal, shoreline re-profiling for spawning, and lake-to-river habitat conversion"
# Row 5
$ws.Range("F5").Value = "The Proponent must conduct post-construction monitoring and submit a report by June 14, 2022, ensuring offsetting measures meet criteria."
$ws.Range("G5").Value = "Structural stability and functionality of offsetting habitat must be maintained; contingency measures required if not met."
$ws.Range("H5").Value = "Prohibition on adverse disturbance of offsetting measures without DFO approval."
$ws.Range("I5").Value = "Proponent must inform DFO of any unauthorized impacts to fish or habitat and comply with all conditions to avoid legal penalties."
$ws.Range("J5").Value = "Riverine, Estuarine, Lentic"
$ws.Range("K5").Value = "None explicitly listed in section 4"
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = "Riparian vegetation (e.g. trees and shrubs and grass)"
$ws.Range("O5").Value = "Habitat enhancement features including woody debris or structures"
$ws.Range("P5").Value = "Pool/deepwater habitat, riffles, mesostructural units"
# Row 6
$ws.Range("F6").Value = "serious harm to fish as described in the authorization"
$ws.Range("G6").Value = "compliance with reporting requirements including annual reports by December 31 each year"
$ws.Range("H6").Value = "implementation of offsetting measures as per the approved plan"
$ws.Range("I6").Value = "provision of access permissions for DFO to monitor offsetting measures"
$ws.Range("L6").Value = "None"
$ws.Range("M6").Value = "Riparian vegetation (e.g. trees and shrubs and grass)"
$ws.Range("N6").Value = "None"
$ws.Range("O6").Value = "None"
$ws.Range("P6").Value = "Pool/deepwater habitat, Riffles (rivers)"
# Row 7
$ws.Range("D7").Value = "53°36'40.96`"N, 108°44'38.01`"W, UTM Coordinates: [Not explicitly converted, but mentioned in the document]"
$ws.Range("F7").Value = "Installation of gravel/boulder bar during low flow period (fall), not extending further than the river water intake structure"
$ws.Range("G7").Value = "Monitoring of offset structure post-construction for three years including sonar surveys and embeddedness surveys"
$ws.Range("H7").Value = "Compliance reports due annually by March 31 from 2022 to 2024"
$ws.Range("I7").Value = "Proponent responsible for design and workmanship; authorization non-transferable"
$ws.Range("L7").Value = "The gravel/boulder bar offset area is specified as requiring a measured area per the Offsetting Plan, but the exact numerical value isn't explicitly stated in the provided text. However, the embeddedness survey and monitoring parameters suggest an area, but without a specific number given."
$ws.Range("M7").Value = "None explicitly mentioned in the context of vegetation cover types like emergent or riparian"
$ws.Range("N7").Value = "Gravel/boulder bar installation as part of offsetting measures"
$ws.Range("O7").Value = "None mentioned"
$ws.Range("P7").Value = "Gravel/boulder bar (as a structure), riffles implied through velocity distribution monitoring"
# Row 8
$ws.Range("F8").Value = "Sedimentation and erosion control measures must be in place and upgraded/maintained to avoid sediment release into the water."
$ws.Range("H8").Value = "All riprap must be clean, free of fine materials, and not obtained from fish-frequented waters below the ordinary high water mark."
$ws.Range("I8").Value = "Water from dewatering must be released into vegetated areas or settling basins and not directly into fish-frequented waters."
$ws.Range("K8").Value = "None explicitly listed in the provided sections"
$ws.Range("L8").Value = "7800 m2"
$ws.Range("M8").Value = "None explicitly mentioned"
$ws.Range("N8").Value = "None explicitly mentioned"
$ws.Range("O8").Value = "None explicitly mentioned"
$ws.Range("P8").Value = "Secondary channel re-graded to increase connectivity during low flow conditions"
# Row 9
$ws.Range("F9").Value = "Inspection reports with dated and annotated photographs during pre-construction, construction, and post-construction periods."
$ws.Range("G9").Value = "Monitoring as per the proposed plan 'Ministry of Transportation – West Region, Highway 401 Grand River Bridge Replacements, Ministry Act Authorization Permit Application Package'."
$ws.Range("H9").Value = "Contingency measures to prevent greater impacts if mitigation measures fail."
$ws.Range("I9").Value = "Completion of offsetting measures by December 31, 2024, with as-built surveys and photographic documentation."
$ws.Range("L9").Value = "150 m2"
$ws.Range("M9").Value = "80% coverage of herbaceous ground cover and seasonal planted stock"
$ws.Range("N9").Value = "50 m2"
$ws.Range("O9").Value = "100 m2"
$ws.Range("P9").Value = "overwintering pool habitat with depths ≥1.7m, boulder clusters as velocity refuge, anchored sweeper trees"
# Row 10
$ws.Range("D10").Value = "Longitude and latitude: 50.894225, -114.009975"
$ws.Range("F10").Value = "measures and standards to avoid and mitigate serious harm to fish shall be implemented before, during, and following in-water work"
$ws.Range("G10").Value = "monitoring and reporting of implementation of offsetting measures with post-construction monitoring for three years"
$ws.Range("H10").Value = "offsetting measures (3,462 m2) to be completed upstream and downstream on the west bank according to the approved plan"
$ws.Range("I10").Value = "contingency measures must be reviewed and approved by DFO if mitigation measures fail or offsetting measures are not effective"
$ws.Range("J10").Value = "riverine"
$ws.Range("L10").Value = "3462 m2"
$ws.Range("M10").Value = "None explicitly mentioned in the offsetting measures"
$ws.Range("P10").Value = "None explicitly mentioned in the offsetting measures"
# Row 11
$ws.Range("G11").Value = "Sediment and erosion control measures must be implemented and monitored, including approved plans and turbidity monitoring."
$ws.Range("H11").Value = "Offsetting measures include placing course rock substrate to provide 2792 HEUs of habitat."
$ws.Range("I11").Value = "Monitoring and reporting requirements for offsetting measures over ten years, including photographic records and annual reports."
$ws.Range("K11").Value = "None"
$ws.Range("L11").Value = "'2792"
$ws.Range("M11").Value = "None"
$ws.Range("N11").Value = "None"
$ws.Range("O11").Value = "None"
$ws.Range("P11").Value = "None"
